$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for new rows 256-269 (column A holds Excel date serial numbers,
# matching the existing data in the sheet)
$data = @(
    @(256, 44330, 1, 17, 66.53880778112647),
    @(257, 44331, 6, 20, 78.28095033073701),
    @(258, 44332, 1, 15, 58.71071274805276),
    @(259, 44333, 4, 16, 62.62476026458961),
    @(260, 44334, 0, 14, 54.79666523151591),
    @(261, 44335, 3, 17, 66.53880778112647),
    @(262, 44336, 2, 17, 66.53880778112647),
    @(263, 44337, 1, 17, 66.53880778112647),
    @(264, 44338, 0, 11, 43.05452268190535),
    @(265, 44339, 1, 11, 43.05452268190535),
    @(266, 44340, 0, 7, 27.39833261575795),
    @(267, 44341, 0, 7, 27.39833261575795),
    @(268, 44342, 1, 5, 19.57023758268425),
    @(269, 44343, 6, 9, 35.22642764883166)
)

# Copy the formatting (number format, font, alignment, border) used by the
# last existing row in column A so the appended date cells match the style
# of the rest of the column.
$ws.Range("A255").Copy() | Out-Null

foreach ($r in $data) {
    $row = $r[0]

    $ws.Range("A$row").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

    $ws.Cells.Item($row, 1).Value = $r[1]
    $ws.Cells.Item($row, 2).Value = $r[2]
    $ws.Cells.Item($row, 3).Value = $r[3]
    $ws.Cells.Item($row, 4).Value = $r[4]
}

$excel.CutCopyMode = 0
